$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Duplicate the "wizard proof" block (rows 22:30) into a fresh block at
# rows 32:40, using a new fixed time (-t0:31) and a new challenge word
# ("WODSE") so the author can resume their first saved game.
# ---------------------------------------------------------------------------

# 1) Clone formatting for the whole block in one shot (reuses existing
#    cellXfs entries instead of synthesizing new styles).
$ws.Range("A22:H30").Copy()
$ws.Range("A32:H40").PasteSpecial(-4122)

# Row 22's block only populates column A - drop the spurious formatted-but
# empty cells the rectangular paste created in columns B:H of row 32, and in
# the unused C/D/G columns of rows 33 and 34 (mirrors rows 23/24 exactly).
$ws.Range("B32:H32").Clear()
$ws.Range("C33:D33").Clear()
$ws.Range("G33").Clear()
$ws.Range("C34:D34").Clear()
$ws.Range("G34").Clear()

# 2) Row 32: new instructions (new fixed time -t0:31)
$ws.Range("A32").Value = "Run Adventure with -d1-JAN-77 -t0:31 to lock the date and time and simplify the calculation."

# 3) Row 33: header labels (identical to rows 13/23)
$ws.Range("A33").Value = "Time"
$ws.Range("B33").Value = "Challenge"
$ws.Range("E33").Value = "MAGNUM"
$ws.Range("F33").Value = "Minutes"
$ws.Range("H33").Value = "Response"

# 4) Row 34: fixed date/time, new challenge word, and the derived formulas
$ws.Range("A34").Value = 28126.021527777779
$ws.Range("B34").Value = "WODSE"
$ws.Range("E34").Value = 11111
$ws.Range("F34").Formula = "=60*HOUR(A34)+MINUTE(A34)"
$ws.Range("H34").Formula = "=CONCAT(H36:H40)"

# 5) Row 35: column headers for the letter-by-letter table (identical to
#    rows 15/25)
$ws.Range("A35").Value = "Y"
$ws.Range("B35").Value = "letter"
$ws.Range("C35").Value = "VAL(Y)"
$ws.Range("D35").Value = "|delta|"
$ws.Range("E35").Value = "D"
$ws.Range("F35").Value = "T"
$ws.Range("G35").Value = "X"
$ws.Range("H35").Value = "letter"

# 6) Rows 36:40: per-letter derivation table, same formulas as rows 26:30
#    but re-pointed at the new block (B34/E34/F34 and the 36:40 range).
$ws.Range("A36").Value = 1
$ws.Range("A37").Value = 2
$ws.Range("A38").Value = 3
$ws.Range("A39").Value = 4
$ws.Range("A40").Value = 5

$ws.Range("B36").Formula = "=UPPER(MID(B34, A36, 1))"
$ws.Range("B37").Formula = "=UPPER(MID(B34, A37, 1))"
$ws.Range("B38").Formula = "=UPPER(MID(B34, A38, 1))"
$ws.Range("B39").Formula = "=UPPER(MID(B34, A39, 1))"
$ws.Range("B40").Formula = "=UPPER(MID(B34, A40, 1))"

$ws.Range("C36").Formula = "=CODE(UPPER(B36)) - CODE(""A"") + 1"
$ws.Range("C37").Formula = "=CODE(UPPER(B37)) - CODE(""A"") + 1"
$ws.Range("C38").Formula = "=CODE(UPPER(B38)) - CODE(""A"") + 1"
$ws.Range("C39").Formula = "=CODE(UPPER(B39)) - CODE(""A"") + 1"
$ws.Range("C40").Formula = "=CODE(UPPER(B40)) - CODE(""A"") + 1"

$ws.Range("D36").Formula = "=ABS(C36-C37)"
$ws.Range("D37").Formula = "=ABS(C37-C38)"
$ws.Range("D38").Formula = "=ABS(C38-C39)"
$ws.Range("D39").Formula = "=ABS(C39-C40)"
$ws.Range("D40").Formula = "=ABS(C40-C36)"

$ws.Range("E36").Formula = "=E34"
$ws.Range("E37").Formula = "=FLOOR(E36/10, 1)"
$ws.Range("E38").Formula = "=FLOOR(E37/10, 1)"
$ws.Range("E39").Formula = "=FLOOR(E38/10, 1)"
$ws.Range("E40").Formula = "=FLOOR(E39/10, 1)"

$ws.Range("F36").Formula = "=40*FLOOR(F34/60, 1)+10*FLOOR(F34/10,1)"
$ws.Range("F37").Formula = "=FLOOR(F36/10, 1)"
$ws.Range("F38").Formula = "=FLOOR(F37/10, 1)"
$ws.Range("F39").Formula = "=FLOOR(F38/10, 1)"
$ws.Range("F40").Formula = "=FLOOR(F39/10, 1)"

$ws.Range("G36").Formula = "=MOD(D36*MOD(E36, 10)+MOD(F36, 10), 26)+1"
$ws.Range("G37").Formula = "=MOD(D37*MOD(E37, 10)+MOD(F37, 10), 26)+1"
$ws.Range("G38").Formula = "=MOD(D38*MOD(E38, 10)+MOD(F38, 10), 26)+1"
$ws.Range("G39").Formula = "=MOD(D39*MOD(E39, 10)+MOD(F39, 10), 26)+1"
$ws.Range("G40").Formula = "=MOD(D40*MOD(E40, 10)+MOD(F40, 10), 26)+1"

$ws.Range("H36").Formula = "=CHAR(64+G36)"
$ws.Range("H37").Formula = "=CHAR(64+G37)"
$ws.Range("H38").Formula = "=CHAR(64+G38)"
$ws.Range("H39").Formula = "=CHAR(64+G39)"
$ws.Range("H40").Formula = "=CHAR(64+G40)"

# 7) Move the active selection to C34, matching where the author's edit
#    left the cursor.
[void]$ws.Range("C34").Select()
